# Rename the "prok only" sheet to "DB peps" (new consolidated sheet name),
# then move the active tab/selection from that sheet back to "ja14_propeps",
# and update the remembered selections on both sheets.

$wb = $excel.ActiveWorkbook

# 1. Rename "prok only" -> "DB peps"
$ws2 = $wb.Worksheets.Item("prok only")
$ws2.Name = "DB peps"

# 2. Update the (now inactive) "DB peps" sheet's remembered selection to A2:A43
$ws2.Range("A2:A43").Select()

# 3. Make "ja14_propeps" the active sheet and set its selection to D15
$ws1 = $wb.Worksheets.Item("ja14_propeps")
$ws1.Activate()
$ws1.Range("D15").Select()
